$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51, shifting existing rows 51-53 down to 52-54
$ws.Rows("51:51").Insert()

# Populate the newly inserted row 51 with the new weekly price entry
$ws.Range("A51").Value = 1
$ws.Range("B51").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C51").Value = "Arica y Parinacota"
$ws.Range("D51").Value = 45166
$ws.Range("E51").Value = 15
$ws.Range("F51").Value = 100112044
$ws.Range("G51").Value = "Perejil"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 100
$ws.Range("K51").Value = 2000
$ws.Range("L51").Value = 2000
$ws.Range("M51").Value = 2000
$ws.Range("N51").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O51").Value = "Región de Arica y Parinacota"
$ws.Range("P51").Value = 1000
$ws.Range("Q51").Value = 2
$ws.Range("R51").Value = "Hortaliza"
